$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.340305
$ws.Range("H2").Value = 1.020915
$ws.Range("I2").Value = 0.2413792532744959
$ws.Range("J2").Value = 0.2522899505114672
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 0.9640756666666667
$ws.Range("N2").Value = 2.892227
$ws.Range("O2").Value = 0.03116597139157128
$ws.Range("P2").Value = 0.03985704515756099
$ws.Range("Q2").Value = 0.328079769745
$ws.Range("R2").Value = 2.952717927705
$ws.Range("S2").Value = 0.007522818902071779
$ws.Range("T2").Value = 0.01005553195033437
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.340305
$ws.Range("H3").Value = 1.020915
$ws.Range("I3").Value = 0.2413792532744959
$ws.Range("J3").Value = 0.2522899505114672
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 6.821209666666667
$ws.Range("N3").Value = 20.463629
$ws.Range("O3").Value = 0.2205113485150815
$ws.Range("P3").Value = 0.2820040699227877
$ws.Range("Q3").Value = 2.321291755615
$ws.Range("R3").Value = 20.891625800535
$ws.Range("S3").Value = 0.05322686464312249
$ws.Range("T3").Value = 0.07114679284485244
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 0.340305
$ws.Range("H4").Value = 1.020915
$ws.Range("I4").Value = 0.2413792532744959
$ws.Range("J4").Value = 0.2522899505114672
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 1.420191666666667
$ws.Range("N4").Value = 4.260575
$ws.Range("O4").Value = 0.04591097398705006
$ws.Range("P4").Value = 0.05871390114682402
$ws.Range("Q4").Value = 0.483298325125
$ws.Range("R4").Value = 4.349684926125001
$ws.Range("S4").Value = 0.01108195661809895
$ws.Range("T4").Value = 0.01481292721466741
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 0.340305
$ws.Range("H5").Value = 1.020915
$ws.Range("I5").Value = 0.2413792532744959
$ws.Range("J5").Value = 0.2522899505114672
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 1.492337666666667
$ws.Range("N5").Value = 4.477012999999999
$ws.Range("O5").Value = 0.04824325997845007
$ws.Range("P5").Value = 0.06169657821187187
$ws.Range("Q5").Value = 0.507849969655
$ws.Range("R5").Value = 4.570649726895
$ws.Range("S5").Value = 0.01164492206912565
$ws.Range("T5").Value = 0.01556542666380002
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 0.340305
$ws.Range("H6").Value = 1.020915
$ws.Range("I6").Value = 0.2413792532744959
$ws.Range("J6").Value = 0.2522899505114672
$ws.Range("K6").Value = 2
$ws.Range("M6").Value = 20.2357845
$ws.Range("N6").Value = 40.471569
$ws.Range("O6").Value = 0.6541684461278472
$ws.Range("P6").Value = 0.5577284055609554
$ws.Range("Q6").Value = 6.886338644272501
$ws.Range("R6").Value = 41.318031865635
$ws.Range("S6").Value = 0.1579026910420771
$ws.Range("T6").Value = 0.1407092718378129
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 0.472982
$ws.Range("H7").Value = 1.418946
$ws.Range("I7").Value = 0.3354874068035369
$ws.Range("J7").Value = 0.3506519309819567
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 0.9640756666666667
$ws.Range("N7").Value = 2.892227
$ws.Range("O7").Value = 0.03116597139157128
$ws.Range("P7").Value = 0.03985704515756099
$ws.Range("Q7").Value = 0.4559904369713333
$ws.Range("R7").Value = 4.103913932742
$ws.Range("S7").Value = 0.01045579092267147
$ws.Range("T7").Value = 0.01397594984773381
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 0.472982
$ws.Range("H8").Value = 1.418946
$ws.Range("I8").Value = 0.3354874068035369
$ws.Range("J8").Value = 0.3506519309819567
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 6.821209666666667
$ws.Range("N8").Value = 20.463629
$ws.Range("O8").Value = 0.2205113485150815
$ws.Range("P8").Value = 0.2820040699227877
$ws.Range("Q8").Value = 3.226309390559334
$ws.Range("R8").Value = 29.036784515034
$ws.Range("S8").Value = 0.07397878048407563
$ws.Range("T8").Value = 0.09888527166319624
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 0.472982
$ws.Range("H9").Value = 1.418946
$ws.Range("I9").Value = 0.3354874068035369
$ws.Range("J9").Value = 0.3506519309819567
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 1.420191666666667
$ws.Range("N9").Value = 4.260575
$ws.Range("O9").Value = 0.04591097398705006
$ws.Range("P9").Value = 0.05871390114682402
$ws.Range("Q9").Value = 0.6717250948833334
$ws.Range("R9").Value = 6.04552585395
$ws.Range("S9").Value = 0.01540255360674006
$ws.Range("T9").Value = 0.02058814281261756
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.472982
$ws.Range("H10").Value = 1.418946
$ws.Range("I10").Value = 0.3354874068035369
$ws.Range("J10").Value = 0.3506519309819567
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 1.492337666666667
$ws.Range("N10").Value = 4.477012999999999
$ws.Range("O10").Value = 0.04824325997845007
$ws.Range("P10").Value = 0.06169657821187187
$ws.Range("Q10").Value = 0.7058488542553333
$ws.Range("R10").Value = 6.352639688298
$ws.Range("S10").Value = 0.01618500618591907
$ws.Range("T10").Value = 0.02163402428497219
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 0.472982
$ws.Range("H11").Value = 1.418946
$ws.Range("I11").Value = 0.3354874068035369
$ws.Range("J11").Value = 0.3506519309819567
$ws.Range("K11").Value = 2
$ws.Range("M11").Value = 20.2357845
$ws.Range("N11").Value = 40.471569
$ws.Range("O11").Value = 0.6541684461278472
$ws.Range("P11").Value = 0.5577284055609554
$ws.Range("Q11").Value = 9.571161824379001
$ws.Range("R11").Value = 57.426970946274
$ws.Range("S11").Value = 0.2194652756041307
$ws.Range("T11").Value = 0.1955685423734369
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 0.4136363333333333
$ws.Range("H12").Value = 1.240909
$ws.Range("I12").Value = 0.2933933655608953
$ws.Range("J12").Value = 0.3066551771687498
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 0.9640756666666667
$ws.Range("N12").Value = 2.892227
$ws.Range("O12").Value = 0.03116597139157128
$ws.Range("P12").Value = 0.03985704515756099
$ws.Range("Q12").Value = 0.3987767238158889
$ws.Range("R12").Value = 3.588990514343
$ws.Range("S12").Value = 0.009143889237547678
$ws.Range("T12").Value = 0.01222236924421473
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 0.4136363333333333
$ws.Range("H13").Value = 1.240909
$ws.Range("I13").Value = 0.2933933655608953
$ws.Range("J13").Value = 0.3066551771687498
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 6.821209666666667
$ws.Range("N13").Value = 20.463629
$ws.Range("O13").Value = 0.2205113485150815
$ws.Range("P13").Value = 0.2820040699227877
$ws.Range("Q13").Value = 2.821500155417889
$ws.Range("R13").Value = 25.393501398761
$ws.Range("S13").Value = 0.06469656668521129
$ws.Range("T13").Value = 0.08647800802448098
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 0.4136363333333333
$ws.Range("H14").Value = 1.240909
$ws.Range("I14").Value = 0.2933933655608953
$ws.Range("J14").Value = 0.3066551771687498
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 1.420191666666667
$ws.Range("N14").Value = 4.260575
$ws.Range("O14").Value = 0.04591097398705006
$ws.Range("P14").Value = 0.05871390114682402
$ws.Range("Q14").Value = 0.5874428736305556
$ws.Range("R14").Value = 5.286985862675
$ws.Range("S14").Value = 0.01346997517423933
$ws.Range("T14").Value = 0.01800492175844778
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 0.4136363333333333
$ws.Range("H15").Value = 1.240909
$ws.Range("I15").Value = 0.2933933655608953
$ws.Range("J15").Value = 0.3066551771687498
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 1.492337666666667
$ws.Range("N15").Value = 4.477012999999999
$ws.Range("O15").Value = 0.04824325997845007
$ws.Range("P15").Value = 0.06169657821187187
$ws.Range("Q15").Value = 0.6172850805352221
$ws.Range("R15").Value = 5.555565724817
$ws.Range("S15").Value = 0.01415425241070671
$ws.Range("T15").Value = 0.0189195751222672
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 0.4136363333333333
$ws.Range("H16").Value = 1.240909
$ws.Range("I16").Value = 0.2933933655608953
$ws.Range("J16").Value = 0.3066551771687498
$ws.Range("K16").Value = 2
$ws.Range("M16").Value = 20.2357845
$ws.Range("N16").Value = 40.471569
$ws.Range("O16").Value = 0.6541684461278472
$ws.Range("P16").Value = 0.5577284055609554
$ws.Range("Q16").Value = 8.370255702703501
$ws.Range("R16").Value = 50.221534216221
$ws.Range("S16").Value = 0.1919286820531903
$ws.Range("T16").Value = 0.1710303030193391
$ws.Range("E17").Value = 2
$ws.Range("G17").Value = 0.182912
$ws.Range("H17").Value = 0.365824
$ws.Range("I17").Value = 0.129739974361072
$ws.Range("J17").Value = 0.09040294133782634
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 0.9640756666666667
$ws.Range("N17").Value = 2.892227
$ws.Range("O17").Value = 0.03116597139157128
$ws.Range("P17").Value = 0.03985704515756099
$ws.Range("Q17").Value = 0.1763410083413333
$ws.Range("R17").Value = 1.058046050048
$ws.Range("S17").Value = 0.00404347232928036
$ws.Range("T17").Value = 0.003603194115278081
$ws.Range("E18").Value = 2
$ws.Range("G18").Value = 0.182912
$ws.Range("H18").Value = 0.365824
$ws.Range("I18").Value = 0.129739974361072
$ws.Range("J18").Value = 0.09040294133782634
$ws.Range("K18").Value = 3
$ws.Range("M18").Value = 6.821209666666667
$ws.Range("N18").Value = 20.463629
$ws.Range("O18").Value = 0.2205113485150815
$ws.Range("P18").Value = 0.2820040699227877
$ws.Range("Q18").Value = 1.247681102549334
$ws.Range("R18").Value = 7.486086615296001
$ws.Range("S18").Value = 0.02860913670267207
$ws.Range("T18").Value = 0.02549399739025805
$ws.Range("E19").Value = 2
$ws.Range("G19").Value = 0.182912
$ws.Range("H19").Value = 0.365824
$ws.Range("I19").Value = 0.129739974361072
$ws.Range("J19").Value = 0.09040294133782634
$ws.Range("K19").Value = 3
$ws.Range("M19").Value = 1.420191666666667
$ws.Range("N19").Value = 4.260575
$ws.Range("O19").Value = 0.04591097398705006
$ws.Range("P19").Value = 0.05871390114682402
$ws.Range("Q19").Value = 0.2597700981333333
$ws.Range("R19").Value = 1.5586205888
$ws.Range("S19").Value = 0.005956488587971717
$ws.Range("T19").Value = 0.005307909361091266
$ws.Range("E20").Value = 2
$ws.Range("G20").Value = 0.182912
$ws.Range("H20").Value = 0.365824
$ws.Range("I20").Value = 0.129739974361072
$ws.Range("J20").Value = 0.09040294133782634
$ws.Range("K20").Value = 3
$ws.Range("M20").Value = 1.492337666666667
$ws.Range("N20").Value = 4.477012999999999
$ws.Range("O20").Value = 0.04824325997845007
$ws.Range("P20").Value = 0.06169657821187187
$ws.Range("Q20").Value = 0.2729664672853334
$ws.Range("R20").Value = 1.637798803712
$ws.Range("S20").Value = 0.006259079312698641
$ws.Range("T20").Value = 0.005577552140832467
$ws.Range("E21").Value = 2
$ws.Range("G21").Value = 0.182912
$ws.Range("H21").Value = 0.365824
$ws.Range("I21").Value = 0.129739974361072
$ws.Range("J21").Value = 0.09040294133782634
$ws.Range("K21").Value = 2
$ws.Range("M21").Value = 20.2357845
$ws.Range("N21").Value = 40.471569
$ws.Range("O21").Value = 0.6541684461278472
$ws.Range("P21").Value = 0.6541684461278472
$ws.Range("Q21").Value = 3.701367814464001
$ws.Range("R21").Value = 14.805471257856
$ws.Range("S21").Value = 0.08487179742844918
$ws.Range("T21").Value = 0.05042028833036646